$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 (VSEARCH), shifting it down to row 12.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the TNT data.
$ws.Range("A11").Value = "TNT"
$ws.Range("B11").Value = 0.4117647058823529
$ws.Range("C11").Value = 0.7167832167832168
$ws.Range("D11").Value = 0.5774647887323944
$ws.Range("E11").Value = 0.8541666666666666
$ws.Range("F11").Value = 0.4361702127659575
